$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text (e.g. thousands separated by
# extra dots, or values like "156.00"/"5.23"). Force those cells to the
# Text number format first so Excel does not reinterpret the replacement
# strings as numbers and strip formatting / introduce rounding.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D23", "D25", "D26", "D29", "D30", "D31", "D32", "D38", "D39", "D44", "D47", "D48", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.637.11"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "2.611.74"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D5").Value = "600.90"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "154.26"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("D9").Value = "2.611.18"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +7.17%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "28.06"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "3.085.76"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "67.609.32"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "2.609.50"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").Value = "365.43"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "2.09"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "69.97"
$ws.Range("E25").Value = "  +4.50%  "
$ws.Range("D26").Value = "10.03"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D29").Value = "579.47"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").Value = "1.02"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "7.90"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").Value = "19.36"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").Value = "155.23"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "41.12"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "156.00"
$ws.Range("D48").Value = "0.0₆0283"
$ws.Range("E48").Value = "  -8.69%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "20.91"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -0.17%  "
